$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 21 (pushes the old "CP Ratio" block and
# everything below it down by 6 rows) in one shot, so literal numbers in
# the rows being pushed down keep their original text representation.
$ws.Rows("21:26").Insert()

# ---- Row 20: abb_4m ----
$ws.Range("A20").Value = "abb_4m"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 4
$ws.Range("D20").Formula = "=VLOOKUP(C20,`$O`$6:`$V`$17,8,TRUE)"
$ws.Range("E20").Formula = "=C20*D20*1000"
$ws.Range("F20").Formula = "=C20*H20*1000"
$ws.Range("G20").Value = 10
$ws.Range("H20").Formula = "=VLOOKUP(C20,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("I20").Formula = "=C20/B20"
$ws.Range("J20").Formula = "=VLOOKUP(I20,`$O`$6:`$S`$18,4,TRUE)"
$ws.Range("K20").Formula = "=C20*J20*1000"
$ws.Range("L20").Formula = "=VLOOKUP(I20,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("M20").Formula = "=L20*C20*1000"

# ---- Row 21: abb_6m ----
$ws.Range("A21").Value = "abb_6m"
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 6.6
$ws.Range("D21").Formula = "=VLOOKUP(C21,`$O`$6:`$V`$17,8,TRUE)"
$ws.Range("E21").Formula = "=C21*D21*1000"
$ws.Range("F21").Formula = "=C21*H21*1000"
$ws.Range("G21").Value = 10
$ws.Range("H21").Formula = "=VLOOKUP(C21,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("I21").Formula = "=C21/B21"
$ws.Range("J21").Formula = "=VLOOKUP(I21,`$O`$6:`$S`$18,4,TRUE)"
$ws.Range("K21").Formula = "=C21*J21*1000"
$ws.Range("L21").Formula = "=VLOOKUP(I21,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("M21").Formula = "=L21*C21*1000"

# ---- Row 22: abb_8m ----
$ws.Range("A22").Value = "abb_8m"
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 8.8
$ws.Range("D22").Formula = "=VLOOKUP(C22,`$O`$6:`$V`$17,8,TRUE)"
$ws.Range("E22").Formula = "=C22*D22*1000"
$ws.Range("F22").Formula = "=C22*H22*1000"
$ws.Range("G22").Value = 10
$ws.Range("H22").Formula = "=VLOOKUP(C22,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("I22").Formula = "=C22/B22"
$ws.Range("J22").Formula = "=VLOOKUP(I22,`$O`$6:`$S`$18,4,TRUE)"
$ws.Range("K22").Formula = "=C22*J22*1000"
$ws.Range("L22").Formula = "=VLOOKUP(I22,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("M22").Formula = "=L22*C22*1000"

# ---- Row 23: abb_4 ----
$ws.Range("A23").Value = "abb_4"
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 4
$ws.Range("D23").Formula = "=VLOOKUP(C23,`$O`$6:`$S`$17,4,TRUE)"
$ws.Range("E23").Formula = "=C23*D23*1000"
$ws.Range("F23").Formula = "=C23*H23*1000"
$ws.Range("G23").Value = 10
$ws.Range("H23").Formula = "=VLOOKUP(C23,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("I23").Formula = "=C23/B23"
$ws.Range("J23").Formula = "=VLOOKUP(I23,`$O`$6:`$S`$18,4,TRUE)"
$ws.Range("K23").Formula = "=C23*J23*1000"
$ws.Range("L23").Formula = "=VLOOKUP(I23,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("M23").Formula = "=L23*C23*1000"

# ---- Row 24: abb_6 ----
$ws.Range("A24").Value = "abb_6"
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 6.6
$ws.Range("D24").Formula = "=VLOOKUP(C24,`$O`$6:`$S`$17,4,TRUE)"
$ws.Range("E24").Formula = "=C24*D24*1000"
$ws.Range("F24").Formula = "=C24*H24*1000"
$ws.Range("G24").Value = 10
$ws.Range("H24").Formula = "=VLOOKUP(C24,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("I24").Formula = "=C24/B24"
$ws.Range("J24").Formula = "=VLOOKUP(I24,`$O`$6:`$S`$18,4,TRUE)"
$ws.Range("K24").Formula = "=C24*J24*1000"
$ws.Range("L24").Formula = "=VLOOKUP(I24,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("M24").Formula = "=L24*C24*1000"

# ---- Row 25: abb_8 ----
$ws.Range("A25").Value = "abb_8"
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 8.8
$ws.Range("D25").Formula = "=VLOOKUP(C25,`$O`$6:`$S`$17,4,TRUE)"
$ws.Range("E25").Formula = "=C25*D25*1000"
$ws.Range("F25").Formula = "=C25*H25*1000"
$ws.Range("G25").Value = 10
$ws.Range("H25").Formula = "=VLOOKUP(C25,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("I25").Formula = "=C25/B25"
$ws.Range("J25").Formula = "=VLOOKUP(I25,`$O`$6:`$S`$18,4,TRUE)"
$ws.Range("K25").Formula = "=C25*J25*1000"
$ws.Range("L25").Formula = "=VLOOKUP(I25,`$O`$6:`$S`$17,5,TRUE)"
$ws.Range("M25").Formula = "=L25*C25*1000"

# Force a full recalc so every new formula gets a fresh cached value
# (writing the formula right after its precedent cell can otherwise
# leave a stale value behind).
$excel.CalculateFull()

# The new rows 20-25 share the row-level "s=1" custom format of the table
# above, matching how Excel auto-extends formatting from the row above
# when new rows are inserted in the middle of a formatted block.
$ws.Range("A20:M25").Style = "Normal"

# Reposition the chart: the twoCellAnchor starts exactly at the inserted
# boundary (row 21) so the engine doesn't auto-shift it during Insert();
# nudge it down by the height of the 6 inserted (default 15pt) rows.
$co = $ws.ChartObjects().Item(1)
$co.Top = $co.Top + 90

# Match the author's final selection.
$ws.Range("A23:H25").Select()
